$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell references to their new text values, derived from the diff.
# These are coin price / percentage-change strings that must stay plain text
# (Excel would otherwise auto-convert numeric-looking values like '17.00' into
# the number 17, dropping the formatting that makes the sheet correct).
$updates = [ordered]@{
    "D2" = '41.757.21'
    "E2" = '  +1.30%  '
    "D3" = '2.262.24'
    "E3" = '  +0.56%  '
    "E4" = '  -0.07%  '
    "D5" = '303.33'
    "E5" = '  +0.46%  '
    "D6" = '92.17'
    "E6" = '  +1.15%  '
    "D7" = '0.530'
    "E7" = '  +1.91%  '
    "E8" = '  -0.08%  '
    "E9" = '  -0.02%  '
    "D10" = '32.49'
    "E10" = '  +2.04%  '
    "D11" = '53.57'
    "E11" = '  -0.53%  '
    "E12" = '  +0.33%  '
    "E13" = '  -1.09%  '
    "D14" = '6.65'
    "E14" = '  +1.13%  '
    "D15" = '2.612.37'
    "E15" = '  +0.58%  '
    "D16" = '14.25'
    "E16" = '  +0.95%  '
    "D17" = '2.277.08'
    "E17" = '  +4.97%  '
    "D18" = '0.771'
    "E18" = '  +2.66%  '
    "D19" = '41.653.29'
    "E19" = '  +1.23%  '
    "D20" = '12.42'
    "E20" = '  +4.20%  '
    "E21" = '  +0.28%  '
    "E22" = '  +1.37%  '
    "D23" = '67.14'
    "E23" = '  +0.46%  '
    "D24" = '239.79'
    "E24" = '  -0.38%  '
    "E25" = '  +0.85%  '
    "E26" = '  +0.01%  '
    "E27" = '  +3.82%  '
    "D28" = '23.94'
    "E28" = '  +0.59%  '
    "D29" = '9.53'
    "E29" = '  -0.24%  '
    "D30" = '2.08'
    "E30" = '  -1.01%  '
    "D31" = '35.36'
    "E31" = '  +6.46%  '
    "D32" = '160.30'
    "E32" = '  +0.75%  '
    "D33" = '5.24'
    "E33" = '  +1.23%  '
    "E34" = '  -0.11%  '
    "D35" = '0.0744'
    "E35" = '  +1.68%  '
    "D36" = '3.02'
    "E36" = '  -0.47%  '
    "D37" = '17.00'
    "E37" = '  +3.18%  '
    "D38" = '2.37'
    "E38" = '  +0.24%  '
    "B39" = 'Kaspa'
    "C39" = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    "D39" = '0.105'
    "E39" = '  +0.53%  '
    "B40" = 'Stellar'
    "C40" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "D40" = '0.116'
    "E40" = '  +1.13%  '
    "E41" = '  +0.41%  '
    "E42" = '  -0.66%  '
    "D43" = '2.013.91'
    "E43" = '  -2.73%  '
    "D44" = '19.29'
    "E44" = '  -4.77%  '
    "E45" = '  +0.96%  '
    "D46" = '10.32'
    "E46" = '  +0.68%  '
    "E47" = '  +5.83%  '
    "E48" = '  -2.45%  '
    "B49" = 'TrustWalletToken'
    "C49" = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    "D49" = '1.15'
    "E49" = '  +1.21%  '
    "B50" = 'Stacks'
    "C50" = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    "D50" = '1.51'
    "E50" = '  +0.50%  '
    "B51" = 'MultiversX'
    "C51" = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    "D51" = '52.39'
    "E51" = '  +3.30%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    # Temporarily force text format so Excel stores the exact string instead of
    # coercing it to a number, then restore the cell's original style so no
    # formatting side-effects remain.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = $origStyle
}

Write-Output ("Updated " + $updates.Count + " cells")
